$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell ref, new value, whether to force Text format
$updates = @(
    @('D2', '26.070.13', $false),
    @('E2', '  +0.43%  ', $false),
    @('D3', '1.639.49', $false),
    @('E3', '  +0.04%  ', $false),
    @('E4', '  +0.45%  ', $false),
    @('D5', '214.76', $true),
    @('E5', '  -0.47%  ', $false),
    @('D6', '0.505', $true),
    @('E6', '  -0.27%  ', $false),
    @('E7', '  +0.42%  ', $false),
    @('E8', '  -1.80%  ', $false),
    @('D9', '0.0626', $true),
    @('E9', '  -1.77%  ', $false),
    @('D10', '18.66', $true),
    @('E10', '  -4.77%  ', $false),
    @('D11', '0.0794', $true),
    @('E11', '  -0.15%  ', $false),
    @('B12', 'Polkadot', $false),
    @('C12', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', $false),
    @('D12', '4.22', $true),
    @('E12', '  -1.43%  ', $false),
    @('B13', 'WrappedEther', $false),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', $false),
    @('D13', '1.630.66', $false),
    @('E13', '  -0.15%  ', $false),
    @('D14', '0.532', $true),
    @('E14', '  -2.08%  ', $false),
    @('D15', '62.41', $true),
    @('E15', '  -0.83%  ', $false),
    @('D16', '0.0₃0749', $false),
    @('D17', '26.074.94', $false),
    @('E17', '  +0.51%  ', $false),
    @('E18', '  +0.46%  ', $false),
    @('D19', '191.21', $true),
    @('E19', '  -0.80%  ', $false),
    @('E20', '  -1.83%  ', $false),
    @('E21', '  -3.24%  ', $false),
    @('D22', '6.16', $true),
    @('E22', '  -1.76%  ', $false),
    @('D23', '144.20', $true),
    @('E23', '  +0.65%  ', $false),
    @('D24', '0.131', $true),
    @('E24', '  -0.85%  ', $false),
    @('E25', '  +0.55%  ', $false),
    @('D26', '1.77', $true),
    @('E26', '  -1.00%  ', $false),
    @('D27', '6.77', $true),
    @('E27', '  -1.61%  ', $false),
    @('D28', '15.26', $true),
    @('E28', '  -2.01%  ', $false),
    @('E29', '  -0.44%  ', $false),
    @('D30', '0.0488', $true),
    @('E30', '  -2.74%  ', $false),
    @('D31', '3.18', $true),
    @('E31', '  -2.15%  ', $false),
    @('E32', '  -3.12%  ', $false),
    @('E33', '  -0.57%  ', $false),
    @('E34', '  -1.10%  ', $false),
    @('E35', '  -2.54%  ', $false),
    @('D36', '1.126.07', $false),
    @('E36', '  -0.62%  ', $false),
    @('E37', '  -0.16%  ', $false),
    @('D38', '0.525', $true),
    @('E38', '  -3.37%  ', $false),
    @('E39', '  -1.40%  ', $false),
    @('D40', '98.82', $true),
    @('E40', '  -0.46%  ', $false),
    @('D41', '0.787', $true),
    @('E41', '  -1.57%  ', $false),
    @('D42', '5.30', $true),
    @('E42', '  -3.21%  ', $false),
    @('D43', '0.0₆0113', $false),
    @('E43', '  -1.95%  ', $false),
    @('D44', '55.32', $true),
    @('E44', '  -2.24%  ', $false),
    @('D45', '0.0527', $true),
    @('E45', '  -0.43%  ', $false),
    @('D46', '1.49', $true),
    @('E46', '  +1.47%  ', $false),
    @('E47', '  -0.04%  ', $false),
    @('D48', '7.63', $true),
    @('E48', '  -0.68%  ', $false),
    @('E49', '  +0.07%  ', $false),
    @('D50', '0.0931', $true),
    @('E50', '  -3.06%  ', $false),
    @('D51', '1.16', $true),
    @('E51', '  -0.60%  ', $false)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    if ($forceText) {
        $ws.Range($ref).NumberFormat = "@"
    }
    $ws.Range($ref).Value = $val
}
